# DividendLiberty / Dividend Calculation.xlsx
# Commit: "updated excel files and added mortgage."
#
# The actual data edit is a correction to December's 401K dividend amount
# on the "Yearly" sheet (cell M14): 97.44 -> 105.36. Every other changed
# cell in the diff is a dependent SUM()/shared formula that recalculates
# automatically as a consequence of that one input changing (O14, M15,
# O15 on "Yearly"; G8, I8, G46, I46 on "All Time", which pull from
# Yearly!M15/O15). The diff also shows the two sheets' last-saved cursor
# position/selection moving, which we replay for fidelity.

$wb = $excel.ActiveWorkbook

# --- "Yearly" sheet -------------------------------------------------
$yearly = $wb.Worksheets.Item("Yearly")

# December row (row 14): correct the 401K dividend amount.
$yearly.Range("M14").Value = 105.36

# The saved selection on this sheet moved from N15 to N14.
$yearly.Range("N14").Select()

# --- "All Time" sheet -------------------------------------------------
$allTime = $wb.Worksheets.Item("All Time")
$allTime.Activate()

# The saved view scrolled so row 19 is at the top, and the selection
# moved from M29 to L19.
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$allTime.Range("L19").Select()

# Recalculate every dependent formula (O14/M15/O15 on Yearly; G8/I8/G46/I46
# on All Time) so the cached <v> values written back to the xlsx match.
$excel.CalculateFull()
